$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 469, shifting the existing rows 469-518 down to 470-519.
$ws.Rows.Item(469).Insert()

# Populate the newly inserted row 469 with the new price-record data.
$ws.Range("A469").Value = 11
$ws.Range("B469").Value = "Vega Monumental Concepción"
$ws.Range("C469").Value = "Bíobío"
$ws.Range("D469").Value = 45132
$ws.Range("E469").Value = 8
$ws.Range("F469").Value = "Fruta"
$ws.Range("G469").Value = 100102
$ws.Range("H469").Value = "Cítricos"
$ws.Range("I469").Value = 100102005
$ws.Range("J469").Value = "Naranja"
$ws.Range("K469").Value = "Navel Late"
$ws.Range("L469").Value = "Primera"
$ws.Range("M469").Value = 350
$ws.Range("N469").Value = 8000
$ws.Range("O469").Value = 8000
$ws.Range("P469").Value = 8000
$ws.Range("Q469").Value = "`$/bandeja 15 kilos empedrada"
$ws.Range("R469").Value = "Región de O'Higgins"
$ws.Range("S469").Value = 533
$ws.Range("T469").Value = 15
